$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Context: rows 548-576 are an 8-rows-per-person block of training-course
# certificates (columns A=Certificate No, B=Name, C=Course name, D=Date,
# E=Result). Two new people are being added, each taking the same set of
# 8 standard courses:
#   - "Ali Khalil Mansour Algallal"   -> rows 548-555 (fill style block #1)
#   - "Mustafa Abraheem Abraheem"     -> rows 556-563 (fill style block #2)
# Certificate numbers DSS1557-DSS1568 are filled down column A for rows
# 558-569; rows 564-569 otherwise stay blank (name/course/date/result to
# be filled in later).
#
# IMPORTANT ORDERING NOTES learned while building this script:
#  * Shared-string indices in the saved file are assigned in first-use
#    order, so the VALUES must be written in the same order the original
#    author entered them for the diff's new shared-string block to line
#    up (Name block1, then certificate numbers 1557-1568, then Name
#    block2).
#  * Strings that look like dates ("05-11-2024" etc.) get auto-converted
#    to real date serials by plain `.Value()` assignment, which also
#    forces a brand-new cell style. Prefixing with a leading apostrophe
#    keeps them as literal text (matching the source file, which stores
#    them as shared strings under a quote-prefixed date-look style).
#  * Because both the apostrophe-text trick and any direct style/
#    NumberFormat change end up inventing a *new* xf record instead of
#    reusing the workbook's existing ones, all formatting is applied
#    LAST via copy / PasteSpecial(xlPasteFormats) from existing same-
#    looking template rows, after every cell's value is already final.
# ---------------------------------------------------------------------------

# --- 1) Write all cell VALUES first, in original-entry order -------------

# 1a. Name for the first new block (rows 548-555)
$ws.Range("B548").Value() = "Ali Khalil Mansour Algallal"

# 1b. Certificate numbers DSS1557-DSS1568 filled down column A, rows 558-569
$certs = @("DSS1557","DSS1558","DSS1559","DSS1560","DSS1561","DSS1562","DSS1563","DSS1564","DSS1565","DSS1566","DSS1567","DSS1568")
for ($i = 0; $i -lt $certs.Length; $i++) {
    $r = 558 + $i
    $ws.Cells.Item($r, 1).Value() = $certs[$i]
}

# 1c. Name for the second new block (rows 556-563)
$ws.Range("B556").Value() = "Mustafa Abraheem Abraheem"

# 1d. Courses / dates / result for both new 8-row blocks
$courses = @("30 Hours Construction Safety & Health","30 Hours G. Industry Safety & Health","Electrical Safety ","Fire Marshal","Scaffold Competent Person","Lifting & Rigging Competent Person","Health & Safety Risk Assessment","Safety Management System & PTW")
$dates = @("05-11-2024","10-11-2024","06-11-2024","03-11-2024","01-11-2024","02-11-2024","07-11-2024","08-11-2024")

for ($i = 0; $i -lt 8; $i++) {
    $r = 548 + $i
    $ws.Cells.Item($r, 2).Value() = "Ali Khalil Mansour Algallal"
    $ws.Cells.Item($r, 3).Value() = $courses[$i]
    $ws.Cells.Item($r, 4).Value() = "'" + $dates[$i]
    $ws.Cells.Item($r, 5).Value() = 1
}

for ($i = 0; $i -lt 8; $i++) {
    $r = 556 + $i
    $ws.Cells.Item($r, 2).Value() = "Mustafa Abraheem Abraheem"
    $ws.Cells.Item($r, 3).Value() = $courses[$i]
    $ws.Cells.Item($r, 4).Value() = "'" + $dates[$i]
    $ws.Cells.Item($r, 5).Value() = 1
}

$ws.Application.CutCopyMode = $false

# --- 2) Apply formatting last, copied from the existing matching blocks --

# 2a. Snapshot formatting rows 548-553 already have for columns A/B
#     (filled cert-no cell + blank placeholder cell) onto rows 564-569,
#     which need that same look even though they stay otherwise blank.
$ws.Range("A548:B553").Copy()
$ws.Range("A564:B569").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 2b. Apply the two alternating fill-color block formats (already used
#     consistently for every 8-row person block elsewhere in the sheet)
#     onto the two new person blocks.
$ws.Range("A532:E539").Copy()
$ws.Range("A548").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A540:E547").Copy()
$ws.Range("A556").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

echo "edit applied"
